$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 2
$ws.Range("A5").Value = 4

# Remove row 6 entirely (shifts nothing below it, just deletes it)
$ws.Rows.Item(6).Delete()
